$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33: Sharman1CON ---
$ws.Range("A33").Value = "Sharman1CON"
$ws.Range("C33").Value = 36647
$ws.Range("D33").Value = 112
$ws.Range("E33").Value = 12.1
$ws.Range("F33").Value = 597
$ws.Range("H33").Value = 237

# --- Row 34: Sharman1CORN ---
$ws.Range("A34").Value = "Sharman1CORN"
$ws.Range("C34").Value = 36647
$ws.Range("D34").Value = 112
$ws.Range("E34").Value = 12.5
$ws.Range("F34").Value = 632
$ws.Range("H34").Value = 229

# --- Row 35: Sharman1LGWP ---
$ws.Range("A35").Value = "Sharman1LGWP"
$ws.Range("C35").Value = 36647
$ws.Range("D35").Value = 138
$ws.Range("E35").Value = 11.8
$ws.Range("F35").Value = 610
$ws.Range("H35").Value = 258

# --- Row 36: Sharman1HGWP ---
$ws.Range("A36").Value = "Sharman1HGWP"
$ws.Range("C36").Value = 36647
$ws.Range("D36").Value = 83
$ws.Range("E36").Value = 11.2
$ws.Range("F36").Value = 588
$ws.Range("H36").Value = 156

# --- Row 37: Sharman2CON ---
$ws.Range("A37").Value = "Sharman2CON"
$ws.Range("C37").Value = 36647
$ws.Range("D37").Value = 112
$ws.Range("E37").Value = 11.9
$ws.Range("F37").Value = 617
$ws.Range("H37").Value = 231

# --- Row 38: Sharman2CORN ---
$ws.Range("A38").Value = "Sharman2CORN"
$ws.Range("C38").Value = 36647
$ws.Range("D38").Value = 146
$ws.Range("E38").Value = 10
$ws.Range("F38").Value = 590
$ws.Range("H38").Value = 223

# --- Row 39: Sharman2LGWP ---
$ws.Range("A39").Value = "Sharman2LGWP"
$ws.Range("C39").Value = 36647
$ws.Range("D39").Value = 128
$ws.Range("E39").Value = 10
$ws.Range("F39").Value = 575
$ws.Range("H39").Value = 197

# --- Row 40: Sharman2HGWP ---
$ws.Range("A40").Value = "Sharman2HGWP"
$ws.Range("C40").Value = 36647
$ws.Range("D40").Value = 131
$ws.Range("E40").Value = 10.5
$ws.Range("F40").Value = 581
$ws.Range("H40").Value = 223

# --- Shared formulas: B33:B36 then G34:G36 then B37:B40 then G38:G40
#     (order matters - it determines the shared-formula group index (si) on save) ---
$ws.Range("B33:B36").Formula = "=C33+D33"
$ws.Range("G33").Formula = "=D33*E33"
$ws.Range("G34:G36").Formula = "=D34*E34"
$ws.Range("B37:B40").Formula = "=C37+D37"
$ws.Range("G37").Formula = "=D37*E37"
$ws.Range("G38:G40").Formula = "=D38*E38"

# --- Copy number formats (date format / 0.000 format) from existing rows ---
$ws.Range("B2:C2").Copy()
$ws.Range("B33:C40").PasteSpecial(-4122)

$ws.Range("F3").Copy()
$ws.Range("F33:F40").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Range("H41").Select()
